# Updates cryptos list values (prices & volume deltas) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.916.41'
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Value = '1.975.79'
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.87%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7121'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.008'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3336'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07007'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8193'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08077'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.29%  '

$ws.Range("D13").Value = '1.983.81'
$ws.Range("E13").Value = '  +0.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.518'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '97.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '266.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.32%  '

$ws.Range("D18").Value = '30.963.38'
$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.034'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000008102'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.62%  '

$ws.Range("D21").Value = '2.241.43'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.006'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.010'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.006'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.841'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.335'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1320'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.03%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.383'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.579'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.563'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.365'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05243'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.37%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.264'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.96%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7713'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.788'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01985'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.879'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.688'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4571'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.072'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8498'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.006'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.981'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.564'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.89'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.556'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4238'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
